# Auto-generated Excel COM-interop script
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values
# on the cryptos worksheet to match the latest scrape, per the GitHub
# Actions commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.401.02"
$ws.Range("D3").Value = "1.846.02"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2935"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "1.860.19"
$ws.Range("E12").Value = "  -6.29%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6803"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001042"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "2.088.53"
$ws.Range("E17").Value = "  -7.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.176"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "29.424.54"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.480"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.351"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.458"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.300"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05650"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.025"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.846"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7103"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.589"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "1.249.98"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01809"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.767"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.381"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9027"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000118"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.088"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.979"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.674"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -0.14%  "
